$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => updated Price (D) / Volume(1h) (E) values scraped for this run.
# Price strings keep their original (locale-style, dot-grouped) text formatting,
# so each is written with a leading "'" to force Excel to store it as text
# instead of auto-coercing it to a number, then the quote-prefix style marker
# is cleared so the cell format stays "General" (matching the source file).
$updates = @(
    @{ Row = 2; D = '28.974.38'; E = '  -4.30%  ' }
    @{ Row = 3; D = '1.956.37'; E = '  -6.73%  ' }
    @{ Row = 4; D = '1.014'; E = '  +1.17%  ' }
    @{ Row = 5; D = '326.69'; E = '  -4.63%  ' }
    @{ Row = 6; D = '1.012'; E = '  +1.00%  ' }
    @{ Row = 7; D = '0.4953'; E = '  -6.15%  ' }
    @{ Row = 8; D = '0.4184'; E = '  -4.49%  ' }
    @{ Row = 9; D = '54.16'; E = '  -1.63%  ' }
    @{ Row = 10; D = '0.08918'; E = '  -4.98%  ' }
    @{ Row = 11; D = '1.092'; E = '  -7.19%  ' }
    @{ Row = 12; D = '22.77'; E = '  -8.07%  ' }
    @{ Row = 13; D = '1.935.12'; E = '  -8.46%  ' }
    @{ Row = 14; D = '7.797'; E = '  -8.88%  ' }
    @{ Row = 15; D = '6.391'; E = '  -6.91%  ' }
    @{ Row = 16; D = '1.014'; E = '  +1.06%  ' }
    @{ Row = 17; D = '0.00001093'; E = '  -5.55%  ' }
    @{ Row = 18; D = '90.39'; E = '  -10.78%  ' }
    @{ Row = 19; D = '0.06663'; E = '  -0.91%  ' }
    @{ Row = 20; D = '19.07'; E = '  -9.62%  ' }
    @{ Row = 21; D = $null; E = '  +0.74%  ' }
    @{ Row = 22; D = '5.919'; E = '  -7.72%  ' }
    @{ Row = 23; D = '29.024.46'; E = '  -4.15%  ' }
    @{ Row = 24; D = '11.86'; E = '  -4.75%  ' }
    @{ Row = 25; D = '2.291'; E = '  -1.28%  ' }
    @{ Row = 26; D = '155.67'; E = '  -4.20%  ' }
    @{ Row = 27; D = '20.50'; E = '  -5.97%  ' }
    @{ Row = 28; D = '6.156'; E = '  -12.22%  ' }
    @{ Row = 29; D = '2.243'; E = '  -10.93%  ' }
    @{ Row = 30; D = '126.14'; E = '  -5.71%  ' }
    @{ Row = 31; D = '1.031'; E = '  -9.10%  ' }
    @{ Row = 32; D = '0.09790'; E = '  -7.03%  ' }
    @{ Row = 33; D = '1.509'; E = '  -9.95%  ' }
    @{ Row = 34; D = '5.757'; E = '  -8.05%  ' }
    @{ Row = 35; D = '3.696'; E = '  -4.62%  ' }
    @{ Row = 36; D = '0.02405'; E = '  -8.49%  ' }
    @{ Row = 37; D = '8.918'; E = '  -11.68%  ' }
    @{ Row = 38; D = '0.06284'; E = '  -7.15%  ' }
    @{ Row = 39; D = '1.279'; E = '  -5.31%  ' }
    @{ Row = 40; D = '0.6390'; E = '  -8.23%  ' }
    @{ Row = 41; D = '11.36'; E = '  -10.46%  ' }
    @{ Row = 42; D = '0.1975'; E = '  -10.71%  ' }
    @{ Row = 43; D = '1.010'; E = '  +0.89%  ' }
    @{ Row = 44; D = '0.6128'; E = '  -9.80%  ' }
    @{ Row = 45; D = '13.33'; E = '  -6.93%  ' }
    @{ Row = 46; D = '2.156'; E = '  -7.43%  ' }
    @{ Row = 47; D = '1.283'; E = '  -1.88%  ' }
    @{ Row = 48; D = '3.470'; E = $null }
    @{ Row = 49; D = '0.00000000328'; E = '  -5.00%  ' }
    @{ Row = 50; D = '0.06851'; E = '  -6.12%  ' }
    @{ Row = 51; D = '1.096'; E = '  -9.62%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Range("D$($u.Row)").Value = "'" + $u.D
        $ws.Range("D$($u.Row)").Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
